# Commit 2019-07-21 kl. 21:57
# Fixzone.xlsx: mark zones 28 and 29 (rows 17-18) as "fixed" like the
# other already-completed rows: copy the green-fill formatting from an
# already-done row onto D:H, and record a 5-point score in column I.
# The dependent formulas in K3/L3/M3 recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "done" look (green fill + thin border) from an already
# completed row (row 9, style index 17 in the OOXML) onto rows 17 & 18.
$doneFormat = $ws.Range("D9:H9")
$doneFormat.Copy()
$ws.Range("D17:H17").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("D18:H18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Record the score for the newly completed zones.
$ws.Range("I17").Value = 5
$ws.Range("I18").Value = 5

# Move the active selection to B28, matching the saved view state.
$ws.Range("B28").Select()
